$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIEfIE")

# Turn on (include) emissions from imported electricity: set control value B2 to 1
$ws.Range("B2").Value = 1

# Reflect the sheet being the active/selected sheet with B3 as the active cell,
# matching the saved view state after the edit.
$ws.Activate()
$ws.Range("B3").Select()
